$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# getAssays / getPrefixes / getClientSites: drop the trailing block
# of empty auto-formatted rows that used to pad the sheets out.
# ---------------------------------------------------------------
$wsGetAssays = $wb.Worksheets.Item("getAssays")
$wsGetAssays.Rows("9:16").Delete()

$wsGetPrefixes = $wb.Worksheets.Item("getPrefixes")
$wsGetPrefixes.Rows("9:12").Delete()

$wsGetClientSites = $wb.Worksheets.Item("getClientSites")
$wsGetClientSites.Rows("9:20").Delete()

# ---------------------------------------------------------------
# putAssays: new "Assert409" dry-run test case block (rows 13-15)
# ---------------------------------------------------------------
$wsPutAssays = $wb.Worksheets.Item("putAssays")

$wsPutAssays.Range("A9:E9").Copy()
$wsPutAssays.Range("A13:E13").PasteSpecial(-4122)

$wsPutAssays.Range("A10:E10").Copy()
$wsPutAssays.Range("A14:E14").PasteSpecial(-4122)

$wsPutAssays.Range("A3:E3").Copy()
$wsPutAssays.Range("A15:E15").PasteSpecial(-4122)

$wsPutAssays.Range("A13").Value = "Assert409"
$wsPutAssays.Range("A14").Value = "EndPoint"
$wsPutAssays.Range("A15").Value = "/configuration/clientCode/assays"
$wsPutAssays.Range("A15:E15").RowHeight = 60

$wsPutAssays.Range("A13:E13").Merge()

$wsPutAssays.Range("D15").Select()

# ---------------------------------------------------------------
# putPrefixes: new "Assert409" dry-run test case block (rows 14-16)
# ---------------------------------------------------------------
$wsPutPrefixes = $wb.Worksheets.Item("putPrefixes")

$wsPutPrefixes.Range("A9:E9").Copy()
$wsPutPrefixes.Range("A14:E14").PasteSpecial(-4122)

$wsPutPrefixes.Range("A10:E10").Copy()
$wsPutPrefixes.Range("A15:E15").PasteSpecial(-4122)

$wsPutPrefixes.Range("A3:E3").Copy()
$wsPutPrefixes.Range("A16:E16").PasteSpecial(-4122)

$wsPutPrefixes.Range("A14").Value = "Assert409"
$wsPutPrefixes.Range("A15").Value = "EndPoint"
$wsPutPrefixes.Range("A16").Value = "/configuration/clientCode/poolingPrefixes"
$wsPutPrefixes.Range("A16:E16").RowHeight = 75

$wsPutPrefixes.Range("A14:E14").Merge()

$wsPutPrefixes.Range("A14:E14").Select()

# ---------------------------------------------------------------
# putClientSites: new "Assert409" dry-run test case block (rows 13-15)
# ---------------------------------------------------------------
$wsPutClientSites = $wb.Worksheets.Item("putClientSites")

$wsPutClientSites.Range("A9:E9").Copy()
$wsPutClientSites.Range("A13:E13").PasteSpecial(-4122)

$wsPutClientSites.Range("A10:E10").Copy()
$wsPutClientSites.Range("A14:E14").PasteSpecial(-4122)

$wsPutClientSites.Range("A3:E3").Copy()
$wsPutClientSites.Range("A15:E15").PasteSpecial(-4122)

$wsPutClientSites.Range("A13").Value = "Assert409"
$wsPutClientSites.Range("A14").Value = "EndPoint"
$wsPutClientSites.Range("A15").Value = "/configuration/clientCode/clientSites"
$wsPutClientSites.Range("A15:E15").RowHeight = 75

$wsPutClientSites.Range("A13:E13").Merge()

$wsPutClientSites.Range("D19").Select()

# ---------------------------------------------------------------
# Active sheet moves from getClientSites to putPrefixes (drop 3).
# ---------------------------------------------------------------
$wsPutPrefixes.Activate()
